# 16.5.1.1a.xlsx - add the 2023 data column (L) to the corruption-perception
# index table on sheet1, matching the upstream "gh-pages" data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New column L: copy column K's formatting (number format / style)
#        down rows 4-14 (header row + the 10 data rows), then overwrite the
#        values with the new 2023 figures.
$ws.Range("K4:K14").Copy()
$ws.Range("L4:L14").PasteSpecial(-4122)   # xlPasteFormats

$ws.Cells.Item(4, 12).Value = 2023

$ws.Cells.Item(5, 12).Value = 22.743990309495757
$ws.Cells.Item(6, 12).Value = 52.401334422687093
$ws.Cells.Item(7, 12).Value = 40.084286291781751
$ws.Cells.Item(8, 12).Value = 58.6564425462321
$ws.Cells.Item(9, 12).Value = 52.689880705632987
$ws.Cells.Item(10, 12).Value = 19.88866894869804
$ws.Cells.Item(11, 12).Value = 35.972443863264772
$ws.Cells.Item(12, 12).Value = 12.061786277026036
$ws.Cells.Item(13, 12).Value = -0.064288010286095529
$ws.Cells.Item(14, 12).Value = 34.132731805770057

# --- 2. Row-height tweaks that came with the refreshed layout.
$ws.Rows.Item(1).RowHeight = 67.5     # title row grew to fit the new column
$ws.Rows.Item(4).RowHeight = 14.25    # year-header row

for ($r = 5; $r -le 14; $r++) {
    $ws.Rows.Item($r).RowHeight = 14.25
}

# --- 3. Reset the lingering UI selection (was parked on M7, outside the
#        table) back to the top-left cell.
$ws.Range("A1").Select()
